$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.028.28"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "3.849.96"
$ws.Range("E3").Value = "  +1.21%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "707.12"
$ws.Range("E5").Value = "  +1.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.21"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").Value = "3.845.56"
$ws.Range("E7").Value = "  +1.12%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.524"
$ws.Range("E9").Value = "  -0.93%  "
$ws.Range("E10").Value = "  -0.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.36"
$ws.Range("E11").Value = "  -1.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.459"
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("E13").Value = "  -1.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.87"
$ws.Range("E14").Value = "  +1.31%  "
$ws.Range("D15").Value = "4.500.75"
$ws.Range("E15").Value = "  +1.22%  "
$ws.Range("D16").Value = "3.824.75"
$ws.Range("E16").Value = "  +0.56%  "
$ws.Range("D17").Value = "71.105.35"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.20"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("E19").Value = "  +0.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.36"
$ws.Range("E20").Value = "  -2.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.69"
$ws.Range("E21").Value = "  -3.89%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "493.68"
$ws.Range("E22").Value = "  +1.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.717"
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.29"
$ws.Range("E24").Value = "  +1.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000147"
$ws.Range("E25").Value = "  +0.94%  "
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.17"
$ws.Range("E26").Value = "  -2.04%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.57"
$ws.Range("E27").Value = "  +0.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.12"
$ws.Range("E28").Value = "  -2.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.19"
$ws.Range("E29").Value = "  +5.18%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.51"
$ws.Range("E31").Value = "  -0.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.27"
$ws.Range("E32").Value = "  -0.58%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.48"
$ws.Range("E33").Value = "  -0.58%  "
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.182"
$ws.Range("E34").Value = "  +0.72%  "
$ws.Range("D35").Value = "3.807.12"
$ws.Range("E35").Value = "  +1.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.16"
$ws.Range("E36").Value = "  -0.95%  "
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.104"
$ws.Range("E38").Value = "  +0.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.38"
$ws.Range("E39").Value = "  +6.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.05"
$ws.Range("E40").Value = "  +1.04%  "
$ws.Range("E41").Value = "  +5.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.33"
$ws.Range("E42").Value = "  -5.66%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "163.32"
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("B46").Value = "FLOKI"
$ws.Range("C46").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.000309"
$ws.Range("E46").Value = "  -6.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "48.81"
$ws.Range("E47").Value = "  -0.94%  "
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "414.85"
$ws.Range("E48").Value = "  +2.90%  "
$ws.Range("B49").Value = "TheGraph"
$ws.Range("C49").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.299"
$ws.Range("E49").Value = "  -0.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.63"
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.23"
$ws.Range("E51").Value = "  -3.97%  "
